$wb = $excel.ActiveWorkbook

$headers = @(
    "Génération",
    "Nombre d'individus départ",
    "Nombre d'individus en vie",
    "Nombre de naissances",
    "Nombre de morts",
    "Moyenne taille individus",
    "Moyenne vue individus",
    "Moyenne vitesse individus",
    "Nombre de morts total"
)

$names = @("simulation1", "simulation2")

foreach ($name in $names) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $name
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }
}
